$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing header (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the numeric data for columns I (I0) and J (IF), rows 2-18
$data = @(
    @(8, 8),
    @(8, 8),
    @(1, 1),
    @(6, 7),
    @(7, 8),
    @(5, 6),
    @(5, 6),
    @(1, 1),
    @(3, 4),
    @(5, 6),
    @(10, 10),
    @(4, 5),
    @(4, 6),
    @(7, 8),
    @(5, 6),
    @(5, 5),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
